# Update Median_Pulse_Width (column C) values per new 8kHz downsampling results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 10
    3 = 28
    4 = 40
    5 = 71
    6 = 9
    7 = 20
    8 = 28
    9 = 49
    10 = 10
    11 = 34
    12 = 55
    13 = 70
    14 = 10
    15 = 29
    16 = 53
    17 = 71
    18 = 21
    19 = 39
    20 = 62
    21 = 80
    23 = 29
    24 = 43
    25 = 64
    26 = 19
    27 = 27
    28 = 34
    29 = 64
    30 = 10
    31 = 28
    32 = 39
    33 = 68
    34 = 9
    35 = 27
    36 = 34
    37 = 62
    38 = 9
    39 = 27
    40 = 61
    41 = 79
    42 = 10
    43 = 28
    44 = 48
    45 = 61
    46 = 12
    47 = 36
    48 = 70
    49 = 74
    50 = 8
    51 = 20
    52 = 28
    53 = 62
    54 = 10
    55 = 26
    56 = 38
    57 = 69
    58 = 27
    59 = 44
    60 = 66
    61 = 70
    62 = 14
    63 = 27
    64 = 34
    65 = 43
    66 = 26
    67 = 48
    68 = 67
    69 = 88
    70 = 18
    71 = 28
    72 = 38
    73 = 58
    74 = 19
    75 = 32
    76 = 39
    77 = 46
    78 = 25
    79 = 54
    80 = 77
    81 = 80
    82 = 15
    83 = 27
    84 = 39
    85 = 74
    86 = 27
    87 = 49
    88 = 66
    89 = 80
    90 = 22
    91 = 30
    92 = 53
    93 = 72
    94 = 19
    95 = 31
    96 = 53
    97 = 71
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
